# Auto update stock data
# Update the "as of" date (cell A2) on every sheet from 2025/10/25 to 2025/10/26.
# The date is stored as a plain text string, so we must avoid Excel's
# auto-detection reinterpreting the literal "2025/10/26" as a real date value.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A2")
    if ($cell.Text -eq "2025/10/25") {
        # Cell stores a plain text date string (not a real date value).
        # Writing a "yyyy/mm/dd"-shaped string directly would make Excel's
        # auto-detection reinterpret it as a date serial, so force a text
        # number format first, then restore the original (General) style
        # so no formatting residue is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = "2025/10/26"
        $cell.Style = "Normal"
    }
}
